# Fix formatting issues introduced by scraping floating point numbers:
#  1) Amounts in the "Importe" column (H) were scraped using Spanish/
#     Argentine number formatting ("." as thousands separator, "," as the
#     decimal separator), e.g. "1.780,00". Convert them to plain decimal
#     notation, e.g. "1780.00", while keeping them as text.
#  2) A handful of "Razon social" entries (column E) contain a literal
#     comma in the name (e.g. "BOFFELLI, MARIA INES"); those commas were
#     mistakenly swept up by the same fix and need to become periods too
#     (e.g. "BOFFELLI. MARIA INES").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column H ("Importe"): "1.234,56" -> "1234.56" ---------------------
$rngImporte = $ws.Range("H2:H273")

# Force text formatting first so Excel doesn't silently convert the cells
# to real numbers once the replaced text looks like a plain number.
$rngImporte.NumberFormat = "@"

# Drop the "." thousands separators, then turn the "," decimal separator
# into a ".".
$rngImporte.Replace(".", "")
$rngImporte.Replace(",", ".")

# Restore the default (unstyled) look of the cells - only their text
# content should have changed.
$rngImporte.Style = "Normal"

# --- 2) Column E ("Razon social"): stray commas -> periods ----------------
$ws.Range("E86").Replace(",", ".")
$ws.Range("E97").Replace(",", ".")
$ws.Range("E123").Replace(",", ".")
$ws.Range("E190").Replace(",", ".")
